$d = $word.ActiveDocument

# The "Class 3 - Extracting and Storing Data" section currently reads (two paragraphs):
#   "Hour 1 - ... Constraints - Lecture Hour 2 - Importing Data into MySQL - Import/Export Wizard - Demographic Data"
#   "Hour 3 - Examples & Exercises - Storing Data for Analysis"
# It needs to become three paragraphs:
#   "Hour 1 - ... Constraints - Lecture "
#   "Hour 2 - Examples & Exercises - Storing Data for Analysis"
#   "Hour 3 - Importing Data into MySQL - Import/Export Wizard - Demographic Data"

# Step 1: pull the "Hour 2 - Importing ... Demographic Data" text out of the Hour 1 paragraph.
# (It gets re-inserted below, re-worded as "Hour 3 - ...", as its own new paragraph.)
$importingRange = $d.Content.Duplicate()
$foundImporting = $importingRange.Find.Execute(
    "Hour 2 - Importing Data into MySQL - Import/Export Wizard – Demographic Data",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundImporting) {
    $importingRange.Text = ""
}

# Step 2: the paragraph that used to read "Hour 3 - Examples & Exercises - Storing Data for Analysis"
# becomes the new "Hour 2" paragraph.
$examplesRange = $d.Content.Duplicate()
$foundExamples = $examplesRange.Find.Execute(
    "Hour 3 - Examples & Exercises - Storing Data for Analysis",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundExamples) {
    $examplesPrefix = $d.Range($examplesRange.Start, $examplesRange.Start + 9)
    $examplesPrefix.Text = "Hour 2 - "

    # Step 3: insert a brand-new paragraph right after it, holding the relocated
    # Importing/Demographic text (renumbered to "Hour 3").
    $breakPos = $examplesRange.End
    $breakRange = $d.Range($breakPos, $breakPos)
    $breakRange.InsertAfter([char]13)

    $newParaRange = $d.Range($breakPos + 1, $breakPos + 1)
    $newParaRange.InsertAfter("Hour 3 - Importing Data into MySQL - Import/Export Wizard – Demographic Data")
}
